# Deploying to gh-pages: update 1.1.1.1a sheet so that the "urban/rural"
# labels become full phrases (городские поселения / сельская местность, etc.)
# and fix the Kyrgyz subtitle wording + the saved sheet view selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 used to hold the short "город" row (urban); expand the labels
$ws.Range("A6").Value = "Шаар жерлери"
$ws.Range("B6").Value = "Городские поселения"
$ws.Range("C6").Value = "City"

# Row 7 used to hold the short "село" row (rural); expand the labels
$ws.Range("A7").Value = "Айыл аймагы"
$ws.Range("B7").Value = "Сельская местность"
$ws.Range("C7").Value = "Village"

# Row 2 subtitle (Kyrgyz) gets extended wording
$ws.Range("A2").Value = "(жалпы калктын санына карата пайыз менен)"

# Restore the saved view: selection on A8 (also resets the scrolled topLeftCell)
$ws.Range("A8").Select()
